# "fetch type lazy loaded"
#
# On slide 15 ("Hibernate - Lazy Loading"), the highlighted phrase
# "lazy loading is by default" inside the content placeholder is
# recoloured from green (00B050) to red (FF0000).
#
# (The many 1/25/2018 -> 2/4/2018 edits visible in the master/layout
# date-placeholder <a:fld type="datetimeFigureOut"> elements are the
# mechanical result of PowerPoint re-stamping the "update
# automatically" date field with the current date whenever the deck
# is opened & saved. That churn isn't a deliberate content edit, so
# it isn't reproduced here.)

$p = $ppt.ActivePresentation
$targetText = "lazy loading is by default"
$maxRunProbe = 50
$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        $textRange = $shape.TextFrame.TextRange

        for ($ri = 1; $ri -le $maxRunProbe; $ri++) {
            $run = $textRange.Runs($ri, 1)
            if ($run.Text -eq $targetText) {
                $red = 255
                $green = 0
                $blue = 0
                $run.Font.Color.RGB = $red + ($green * 256) + ($blue * 65536)
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not locate run '$targetText' to recolor."
}
